$d = $word.ActiveDocument

$d.Content.Find.Execute("58×70=", $true, $false, $false, $false, $false, $true, 1, $false, "60×24=", 2) | Out-Null
$d.Content.Find.Execute("60×63=", $true, $false, $false, $false, $false, $true, 1, $false, "35×73=", 2) | Out-Null
$d.Content.Find.Execute("94×45=", $true, $false, $false, $false, $false, $true, 1, $false, "33×90=", 2) | Out-Null
$d.Content.Find.Execute("60×73=", $true, $false, $false, $false, $false, $true, 1, $false, "26×85=", 2) | Out-Null
$d.Content.Find.Execute("71×83=", $true, $false, $false, $false, $false, $true, 1, $false, "42×77=", 2) | Out-Null
$d.Content.Find.Execute("39×16=", $true, $false, $false, $false, $false, $true, 1, $false, "60×14=", 2) | Out-Null
$d.Content.Find.Execute("17×97=", $true, $false, $false, $false, $false, $true, 1, $false, "96×63=", 2) | Out-Null
$d.Content.Find.Execute("97×48=", $true, $false, $false, $false, $false, $true, 1, $false, "66×53=", 2) | Out-Null
$d.Content.Find.Execute("46×24=", $true, $false, $false, $false, $false, $true, 1, $false, "70×37=", 2) | Out-Null
$d.Content.Find.Execute("50×44=", $true, $false, $false, $false, $false, $true, 1, $false, "20×33=", 2) | Out-Null
$d.Content.Find.Execute("66×40=", $true, $false, $false, $false, $false, $true, 1, $false, "72×70=", 2) | Out-Null
$d.Content.Find.Execute("80×69=", $true, $false, $false, $false, $false, $true, 1, $false, "17×63=", 2) | Out-Null
$d.Content.Find.Execute("20×58=", $true, $false, $false, $false, $false, $true, 1, $false, "32×15=", 2) | Out-Null
$d.Content.Find.Execute("76×97=", $true, $false, $false, $false, $false, $true, 1, $false, "91×65=", 2) | Out-Null
$d.Content.Find.Execute("64×26=", $true, $false, $false, $false, $false, $true, 1, $false, "24×95=", 2) | Out-Null
$d.Content.Find.Execute("94×49=", $true, $false, $false, $false, $false, $true, 1, $false, "98×67=", 2) | Out-Null
$d.Content.Find.Execute("38×56=", $true, $false, $false, $false, $false, $true, 1, $false, "70×30=", 2) | Out-Null
$d.Content.Find.Execute("73×99=", $true, $false, $false, $false, $false, $true, 1, $false, "48×79=", 2) | Out-Null
$d.Content.Find.Execute("80×61=", $true, $false, $false, $false, $false, $true, 1, $false, "62×95=", 2) | Out-Null
$d.Content.Find.Execute("26×55=", $true, $false, $false, $false, $false, $true, 1, $false, "11×14=", 2) | Out-Null
$d.Content.Find.Execute("57×50=", $true, $false, $false, $false, $false, $true, 1, $false, "78×20=", 2) | Out-Null
$d.Content.Find.Execute("19×91=", $true, $false, $false, $false, $false, $true, 1, $false, "49×30=", 2) | Out-Null
$d.Content.Find.Execute("93×65=", $true, $false, $false, $false, $false, $true, 1, $false, "11×81=", 2) | Out-Null
$d.Content.Find.Execute("62×19=", $true, $false, $false, $false, $false, $true, 1, $false, "55×54=", 2) | Out-Null
$d.Content.Find.Execute("62×98=", $true, $false, $false, $false, $false, $true, 1, $false, "87×58=", 2) | Out-Null
